# Fix typo in pptx: slide 22, "reservation_time" field description
# was prompting for a birthday ("생일을 입력하세요") instead of a
# time ("시간을 입력하세요"). Also add a small leftover empty
# textbox shape that appeared on the same slide in the source commit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(22)

# --- 1. Fix the typo: "생일" (birthday) -> "시간" (time) -------------
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$fullText = $tr.Text
$idx = $fullText.IndexOf("생일")
if ($idx -ge 0) {
    $bad = $tr.Characters($idx + 1, 2)
    $bad.Text = "시간"
}

# --- 2. Add the new (empty) textbox shape ----------------------------
$newBox = $s.Shapes.AddTextbox(1, 562.6285826771654, 305.4857480314961, 14.545748031496062, 29.081259842519685)
$newBox.Fill.Visible = 0
$newBox.TextFrame.WordWrap = 0
$newBox.TextFrame.AutoSize = 1
